# Updated mean summary results
#
# The underlying data didn't change, but the rows were re-derived: each
# (species, season) group's "unlabelled" habitat_type row (blank C cell -
# the NA/overall bucket) was dropped, and the season block order changed
# from Fall, Spring, Summer, Winter -> Fall, Winter, Spring, Summer.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- read the existing data block (rows 2..47, cols A..E) -----------------
$lastRow = $ws.UsedRange.Rows.Count
$data = $ws.Range("A2:E$lastRow").Value2
$rowCount = $data.GetLength(0)

$seasonOrder = @{ "Fall" = 0; "Winter" = 1; "Spring" = 2; "Summer" = 3 }

$fishOrder = @{}
$nextFishRank = 0

$records = @()
for ($i = 1; $i -le $rowCount; $i++) {
    $species = $data[$i, 1]
    $season  = $data[$i, 2]
    $habitat = $data[$i, 3]
    $accel   = $data[$i, 4]
    $sem     = $data[$i, 5]

    # drop the blank-habitat ("NA" bucket) rows entirely
    if ($null -eq $habitat -or $habitat -eq "") {
        continue
    }

    if (-not $fishOrder.ContainsKey($species)) {
        $fishOrder[$species] = $nextFishRank
        $nextFishRank++
    }

    $records += [PSCustomObject]@{
        Species    = $species
        Season     = $season
        Habitat    = $habitat
        Accel      = $accel
        Sem        = $sem
        FishRank   = $fishOrder[$species]
        SeasonRank = $seasonOrder[$season]
    }
}

$sorted = $records | Sort-Object FishRank, SeasonRank

# --- clear the old data rows below the header ------------------------------
$ws.Range("A2:E$lastRow").ClearContents()

# --- write the reordered/filtered rows back out ----------------------------
$r = 2
foreach ($rec in $sorted) {
    $ws.Cells.Item($r, 1).Value = $rec.Species
    $ws.Cells.Item($r, 2).Value = $rec.Season
    $ws.Cells.Item($r, 3).Value = $rec.Habitat
    $ws.Cells.Item($r, 4).Value = $rec.Accel
    $ws.Cells.Item($r, 5).Value = $rec.Sem
    $r++
}

Write-Output "Rewrote $($sorted.Count) data rows (from $rowCount)."
